$wb = $excel.ActiveWorkbook

# --- Remove negative signs from the production and degradation rate values ---
$ws1 = $wb.Worksheets.Item("production_rates")
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws1.Range("B$r")
    $cell.Value2 = -1 * $cell.Value2
}

$ws2 = $wb.Worksheets.Item("degradation_rates")
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws2.Range("B$r")
    $cell.Value2 = -1 * $cell.Value2
}

# --- Update selections / active sheet to match the saved view state ---
$ws1.Activate()
$ws1.Range("C1:C1048576").Select()

$ws2.Activate()
$ws2.Range("C1:C1048576").Select()
